$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.429735422134399
$ws.Range("B1").Value = 3.516790390014648
$ws.Range("C1").Value = 5.37604284286499
$ws.Range("D1").Value = 1.736026763916016
$ws.Range("E1").Value = 0.9735450744628906
